$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is a plain text or a multi-dot string that Excel will not
# auto-convert to a number (safe to assign directly). ---
$ws.Range('D2').Value = '43.662.74'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '2.331.37'
$ws.Range('E3').Value = '  +4.55%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +6.57%  '
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('E11').Value = '  +2.86%  '
$ws.Range('E12').Value = '  +2.62%  '
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('D14').Value = '2.678.12'
$ws.Range('E14').Value = '  +4.43%  '
$ws.Range('E15').Value = '  +3.67%  '
$ws.Range('E16').Value = '  +8.60%  '
$ws.Range('D17').Value = '2.325.92'
$ws.Range('E17').Value = '  +4.45%  '
$ws.Range('D18').Value = '43.654.48'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('E19').Value = '  +4.98%  '
$ws.Range('E20').Value = '  +8.31%  '
$ws.Range('E21').Value = '  +2.61%  '
$ws.Range('E22').Value = '  +3.52%  '
$ws.Range('E23').Value = '  -2.87%  '
$ws.Range('E24').Value = '  +8.37%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('E26').Value = '  +1.92%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E27').Value = '  +3.52%  '
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E30').Value = '  +9.52%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('E31').Value = '  -1.47%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('E33').Value = '  -1.36%  '
$ws.Range('E34').Value = '  +2.38%  '
$ws.Range('E36').Value = '  +3.53%  '
$ws.Range('E37').Value = '  -2.51%  '
$ws.Range('E38').Value = '  +3.46%  '
$ws.Range('E39').Value = '  -1.81%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E40').Value = '  +9.23%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E41').Value = '  +10.10%  '
$ws.Range('E42').Value = '  +19.56%  '
$ws.Range('E43').Value = '  -3.26%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E45').Value = '  +9.15%  '
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('E49').Value = '  +2.19%  '
$ws.Range('E50').Value = '  +17.68%  '
$ws.Range('D51').Value = '2.557.86'
$ws.Range('E51').Value = '  +4.41%  '

# --- Price cells whose new value looks like a plain decimal number. Excel would silently
# coerce a bare numeric string into a Double (dropping trailing zeros / using scientific
# notation), so we force the cell to Text format first, assign the literal digits, then
# restore the default "Normal" style so no stray formatting is left behind. ---
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '271.12'
$c.Style = 'Normal'

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '95.43'
$c.Style = 'Normal'

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '45.20'
$c.Style = 'Normal'

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '8.12'
$c.Style = 'Normal'

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '15.65'
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '0.861'
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.0000108'
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '6.46'
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '72.18'
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '239.91'
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '2.28'
$c.Style = 'Normal'

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '9.40'
$c.Style = 'Normal'

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.54'
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '11.45'
$c.Style = 'Normal'

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '3.48'
$c.Style = 'Normal'

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.28'
$c.Style = 'Normal'

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '22.70'
$c.Style = 'Normal'

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '38.21'
$c.Style = 'Normal'

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '172.49'
$c.Style = 'Normal'

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.0901'
$c.Style = 'Normal'

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.49'
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '3.37'
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.35'
$c.Style = 'Normal'

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.235'
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.37'
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '12.12'
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '62.27'
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '9.16'
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '5.36'
$c.Style = 'Normal'

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '100.49'
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.190'
$c.Style = 'Normal'

Write-Host "Applied cryptos update"